# Trade #28 closed at 2026-02-17 08:03:36 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Summary" sheet - update headline stats after the new closed trade
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.39   # Current Capital
$summary.Range("B4").Value = -0.61     # Total P&L $
$summary.Range("B5").Value = -0.44     # Total P&L %
$summary.Range("B6").Value = 28        # Total Trades
$summary.Range("B8").Value = 14        # Losing Trades
$summary.Range("B9").Value = 25        # Win Rate %

# ---------------------------------------------------------------------------
# 2. "Strategy Status" sheet - MarketMaking row (row 4) stats
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.39
$status.Range("D4").Value = 28
$status.Range("E4").Value = -0.61
$status.Range("F4").Value = -0.61
$status.Range("G4").Value = 25

# ---------------------------------------------------------------------------
# 3/4. "All Trades" and "MarketMaking" sheets - append new trade #28 as row 29
# ---------------------------------------------------------------------------
$tradeSheets = @($wb.Worksheets.Item("All Trades"), $wb.Worksheets.Item("MarketMaking"))

foreach ($ws in $tradeSheets) {
    $ws.Range("A29").Value = 28

    # Columns B and C hold plain text (date / time) in this workbook, not
    # real Excel dates - force text storage so Excel doesn't auto-convert
    # the string to a date serial, then drop the formatting change so the
    # cell keeps using the sheet's default (unstyled) look.
    $ws.Range("B29").NumberFormat = "@"
    $ws.Range("B29").Value = "2026-02-17"
    $ws.Range("C29").NumberFormat = "@"
    $ws.Range("C29").Value = "08:03:30"
    $ws.Range("B29:C29").ClearFormats()

    $ws.Range("D29").Value = "MarketMaking"
    $ws.Range("E29").Value = "UP"
    $ws.Range("F29").Value = 0.33
    $ws.Range("G29").Value = 0.3
    $ws.Range("H29").Value = "CLOSED"
    $ws.Range("I29").Value = -9.0909
    $ws.Range("J29").Value = -0.03
    $ws.Range("K29").Value = 99.39
    $ws.Range("L29").Value = 0
    $ws.Range("M29").Value = 0
    $ws.Range("N29").Value = 0.6
    $ws.Range("O29").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P29").Value = "early_exit"
    $ws.Range("Q29").Value = 0.13
}
